# Adds the "logout" Q&A entries (3 new rows of chat transcript) plus the
# small logout/user icon pair that sits to the left of the new answer row,
# and nudges the sheet view/selection to where the new content now is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New transcript rows 37-39 (question/answer/question), and the
#    trailing numbered-but-empty rows 40-45 that extend column A.
#    Rows 35-36 also pick up their column-A running numbers here.
# ---------------------------------------------------------------------

$ws.Range("A35").Value = 34
$ws.Range("A36").Value = 35

# Row 37: answer row -> gets the "AI answer" look (Segoe UI 12, font
# color 343541, vertically centered + wrapped), matching the style
# already used for every other answer cell in the sheet.
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "Hogyan kell kitörölni a böngéső localStorage tartalmát JavaScript használatával"

$ans37 = $ws.Range("B37")
$ans37.WrapText = $true
$ans37.VerticalAlignment = -4108
$ans37.Font.Name = "Segoe UI"
$ans37.Font.Size = 12
$ans37.Font.Color = 4273460
$ws.Rows.Item(37).RowHeight = 17.25

# Row 38: plain question row (default style).
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = 'Hogyan kell beletenni a tokent az alábbi hívásba? (Fejrészben [{"key":"Authorization","value":"Token eeb1dca2fb50531e82cb8433aa458d3d3d8abc52b8c235d98ca561a90b192662","type":"text","enabled":false}]'

# Row 39: plain question row (default style).
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "Hogyan kell a főoldalt betölteni és friisíteni is? (A navbár különben nem frissül)"

# Rows 40-45: only column A keeps counting, column B left blank.
$ws.Range("A40").Value = 39
$ws.Range("A41").Value = 40
$ws.Range("A42").Value = 41
$ws.Range("A43").Value = 42
$ws.Range("A44").Value = 43
$ws.Range("A45").Value = 44

# ---------------------------------------------------------------------
# 2) Small logout/user icon pair, anchored near row 36 in column B.
# ---------------------------------------------------------------------

$icon1 = $ws.Shapes.AddShape(1, 48, 595.5, 24, 24)
$icon1.Name = "AutoShape 1"

$icon2 = $ws.Shapes.AddShape(1, 72.75, 595.5, 24, 24)
$icon2.Name = "AutoShape 2"

# ---------------------------------------------------------------------
# 3) View state: scroll to the top of the sheet and select the new
#    last question cell.
# ---------------------------------------------------------------------

try {
    $excel.ActiveWindow.ScrollRow = 19
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}

$ws.Range("B39").Select() | Out-Null
